$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 currently holds "mngr137319"/"bAdYvyb" — update it to duplicate
# the data-driven credentials already used in row 2 ("mngr255784"/"esYnezY"),
# as part of the new data-driven TestNG test case.
$ws.Range("A4").Value = "mngr255784"
$ws.Range("B4").Value = "esYnezY"

# Reflect the newly active/selected cell in the sheet view.
$ws.Range("A4:B4").Select()

$wb.Save()
